$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (hour 1)
$ws.Range("B2").Value = 12500
$ws.Range("I2").Value = 10000
$ws.Range("P2").Value = 0
$ws.Range("S2").Value = 65724
$ws.Range("T2").Value = 5632.2595
$ws.Range("U2").Value = 22500

# Row 3 (hour 2)
$ws.Range("B3").Value = 12500
$ws.Range("I3").Value = 10000
$ws.Range("P3").Value = 0
$ws.Range("S3").Value = 63085
$ws.Range("T3").Value = 5478.844000000001
$ws.Range("U3").Value = 22500
$ws.Range("V3").ClearContents()
$ws.Range("W3").ClearContents()

# Row 4 (hour 3)
$ws.Range("B4").Value = 12500
$ws.Range("I4").Value = 10000
$ws.Range("S4").Value = 61258
$ws.Range("T4").Value = 5431.737499999999
$ws.Range("U4").Value = 22500

# Row 5 (hour 4)
$ws.Range("B5").Value = 12500
$ws.Range("I5").Value = 10000
$ws.Range("S5").Value = 60272
$ws.Range("T5").Value = 5351.2935
$ws.Range("U5").Value = 22500

# Row 6 (hour 5)
$ws.Range("I6").Value = 10000
$ws.Range("S6").Value = 63182
$ws.Range("T6").Value = 5443.7565
$ws.Range("U6").Value = 22500

# Row 7 (hour 6)
$ws.Range("S7").Value = 67802
$ws.Range("T7").Value = 5521.645500000001

# Row 8 (hour 7)
$ws.Range("S8").Value = 69571
$ws.Range("T8").Value = 6568.996000000001

# Row 9 (hour 8)
$ws.Range("S9").Value = 81158
$ws.Range("T9").Value = 7039.0705

# Row 10 (hour 9)
$ws.Range("B10").Value = 25000
$ws.Range("I10").Value = 20000
$ws.Range("P10").Value = 12000
$ws.Range("S10").Value = 98453
$ws.Range("T10").Value = 8768.984
$ws.Range("U10").Value = 57000

# Row 11 (hour 10)
$ws.Range("B11").Value = 25000
$ws.Range("P11").Value = 20000
$ws.Range("S11").Value = 108429
$ws.Range("T11").Value = 14287.7385
$ws.Range("U11").Value = 65000

# Row 12 (hour 11)
$ws.Range("P12").Value = 20000
$ws.Range("S12").Value = 111739
$ws.Range("T12").Value = 16203.32
$ws.Range("U12").Value = 65000

# Row 13 (hour 12)
$ws.Range("S13").Value = 113097
$ws.Range("T13").Value = 16253.0025

# Row 14 (hour 13)
$ws.Range("S14").Value = 112752
$ws.Range("T14").Value = 16193.9575

# Row 15 (hour 14)
$ws.Range("S15").Value = 116975
$ws.Range("T15").Value = 16307.025

# Row 16 (hour 15)
$ws.Range("S16").Value = 117642
$ws.Range("T16").Value = 16253.9055

# Row 17 (hour 16)
$ws.Range("S17").Value = 96422
$ws.Range("T17").Value = 15667.071

# Row 18 (hour 17)
$ws.Range("S18").Value = 94648
$ws.Range("T18").Value = 16597.392
$ws.Range("V18").Value = 4270.840211111111
$ws.Range("W18").Value = 5.076102251596674

# Row 19 (hour 18)
$ws.Range("S19").ClearContents()
$ws.Range("T19").Value = 16121.5775

# Row 20 (hour 19)
$ws.Range("S20").ClearContents()
$ws.Range("T20").Value = 15277.878

# Row 21 (hour 20)
$ws.Range("S21").ClearContents()
$ws.Range("T21").Value = 13983.214

# Row 22 (hour 21)
$ws.Range("S22").ClearContents()
$ws.Range("T22").Value = 12008.9865

# Row 23 (hour 22)
$ws.Range("S23").ClearContents()
$ws.Range("T23").Value = 9648.002000000002

# Row 24 (hour 23)
$ws.Range("S24").ClearContents()
$ws.Range("T24").Value = 6756.816500000001

# Row 25 (hour 24)
$ws.Range("S25").ClearContents()
$ws.Range("T25").Value = 6005.93
